$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append two new weekly rows ---
$wsData = $wb.Worksheets.Item("Data")

# Copy the style used by the last existing date cell (A109) onto the new
# date cells so they keep the same bordered/centered/date-numfmt style.
$wsData.Range("A109").Copy()
$wsData.Range("A110:A111").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsData.Cells.Item(110, 1).Value = 45231
$wsData.Cells.Item(110, 2).Value = 820.487

$wsData.Cells.Item(111, 1).Value = 45238
$wsData.Cells.Item(111, 2).Value = 772.686

# --- Sheet "SeriesInfo": update metadata values ---
# These cells hold date-looking text (FRED metadata strings), not real
# dates. Briefly force a text format so Excel doesn't auto-convert the
# string into a date serial, then restore the cell to the default
# "Normal" style (matching the original, un-styled inline-string cells).
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Cells.Item(3, 2).Value = "2023-11-15"
$wsInfo.Range("B3").Style = "Normal"

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Cells.Item(4, 2).Value = "2023-11-15"
$wsInfo.Range("B4").Style = "Normal"

$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Cells.Item(7, 2).Value = "2023-11-08"
$wsInfo.Range("B7").Style = "Normal"

$wsInfo.Range("B14").NumberFormat = "@"
$wsInfo.Cells.Item(14, 2).Value = "2023-11-09 15:38:01-06"
$wsInfo.Range("B14").Style = "Normal"
